$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell writes. Numeric-looking text values are forced to stay
# text (matching the source inlineStr cells) by temporarily applying a Text
# number format, then clearing formatting again so no visible style changes.

$ws.Range("D2").Value = "34.071.43"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "1.785.43"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0705"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  -3.16%  "
$ws.Range("D13").Value = "1.776.39"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "34.055.65"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.63"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.76"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.75"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -4.70%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.66"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E34").Value = "  -7.54%  "
$ws.Range("D35").Value = "1.394.53"
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.645"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.20"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.913"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.21%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "78.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("D44").Value = "0.0₆0143"
$ws.Range("E44").Value = "  +11.12%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.27"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "1.941.45"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("E51").Value = "  -0.33%  "
